# The workbook keeps a weekly price log for "Ciboulette" at Vega Central
# Mapocho de Santiago. This commit adds one new weekly record, inserted as
# row 203 (pushing every existing record from row 203 down by one row, all
# the way to the end of the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 203; everything below shifts down.
$ws.Rows(203).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A203").Value = 9
$ws.Range("B203").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C203").Value = "Metropolitana"
$ws.Range("D203").Value = 44523
$ws.Range("E203").Value = 13
$ws.Range("F203").Value = 100112039
$ws.Range("G203").Value = "Ciboulette"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 250
$ws.Range("K203").Value = 800
$ws.Range("L203").Value = 1000
$ws.Range("M203").Value = 900
$ws.Range("N203").Value = "`$/docena de atados"
$ws.Range("O203").Value = "Región Metropolitana"
$ws.Range("P203").Value = 300
$ws.Range("Q203").Value = 3
$ws.Range("R203").Value = "Hortaliza"
